$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("G8").Value = 3.6
$ws.Range("H8").Value = 3
$ws.Range("I8").Value = 2.25
$ws.Range("L8").Value = 1.57
$ws.Range("M8").Value = 2.25
$ws.Range("N8").Value = 2.88
$ws.Range("O8").Value = 1.4
$ws.Range("T8").Value = 7.5
$ws.Range("X8").Value = 41
$ws.Range("AF8").Value = 9

# Row 9
$ws.Range("G9").Value = 1.67
$ws.Range("H9").Value = 3.5
$ws.Range("I9").Value = 6.25
$ws.Range("P9").Value = 1.57
$ws.Range("Q9").Value = 2.25
$ws.Range("R9").Value = 2.38
$ws.Range("S9").Value = 1.53
$ws.Range("V9").Value = 9.5
$ws.Range("W9").Value = 12
$ws.Range("AB9").Value = 23
$ws.Range("AC9").Value = 101
$ws.Range("AF9").Value = 29
$ws.Range("AG9").Value = 21
$ws.Range("AJ9").Value = 67

# Row 12
$ws.Range("G12").Value = 1.85
$ws.Range("I12").Value = 4.5
$ws.Range("L12").Value = 1.62
$ws.Range("M12").Value = 2.2
$ws.Range("N12").Value = 2.88
$ws.Range("O12").Value = 1.4
$ws.Range("T12").Value = 4.5
$ws.Range("AB12").Value = 29
$ws.Range("AE12").Value = 8

# Row 14
$ws.Range("G14").Value = 1.95
$ws.Range("H14").Value = 2.95
$ws.Range("I14").Value = 4.4
$ws.Range("J14").Value = 1.09
$ws.Range("K14").Value = 6
$ws.Range("L14").Value = 1.39
$ws.Range("M14").Value = 2.77
$ws.Range("Q14").Value = 2.52
$ws.Range("R14").Value = 1.87
$ws.Range("S14").Value = 1.83
$ws.Range("T14").Value = 6.1
$ws.Range("U14").Value = 8.75
$ws.Range("V14").Value = 8.25
$ws.Range("X14").Value = 17
$ws.Range("Y14").Value = 30
$ws.Range("Z14").Value = 6
$ws.Range("AA14").Value = 5.7
$ws.Range("AB14").Value = 14.5
$ws.Range("AD14").Value = 700
$ws.Range("AE14").Value = 10.5
$ws.Range("AF14").Value = 25
$ws.Range("AG14").Value = 14
$ws.Range("AH14").Value = 80
$ws.Range("AI14").Value = 45
$ws.Range("AJ14").Value = 50

# Row 19
$ws.Range("G19").Value = 2.2
$ws.Range("H19").Value = 3.25
$ws.Range("I19").Value = 3.2
$ws.Range("R19").Value = 1.62
$ws.Range("S19").Value = 2.2
$ws.Range("T19").Value = 9.5
$ws.Range("U19").Value = 12
$ws.Range("V19").Value = 9.5
$ws.Range("W19").Value = 21
$ws.Range("X19").Value = 17
$ws.Range("Y19").Value = 23
$ws.Range("AE19").Value = 12
$ws.Range("AF19").Value = 17
$ws.Range("AG19").Value = 12
$ws.Range("AH19").Value = 34
$ws.Range("AI19").Value = 23
$ws.Range("AJ19").Value = 29

# Row 20
$ws.Range("G20").Value = 2.75
$ws.Range("I20").Value = 2.55
$ws.Range("N20").Value = 2.3
$ws.Range("O20").Value = 1.6
$ws.Range("U20").Value = 13
$ws.Range("V20").Value = 11
$ws.Range("W20").Value = 29
$ws.Range("Y20").Value = 41
$ws.Range("AA20").Value = 6
$ws.Range("AG20").Value = 10

# Row 25
$ws.Range("G25").Value = 3.6
$ws.Range("H25").Value = 3.25
$ws.Range("I25").Value = 2
$ws.Range("N25").Value = 2.15
$ws.Range("O25").Value = 1.67
$ws.Range("R25").Value = 1.91
$ws.Range("S25").Value = 1.8
$ws.Range("X25").Value = 34
$ws.Range("AB25").Value = 17
$ws.Range("AD25").Value = 351
$ws.Range("AE25").Value = 6.5
$ws.Range("AF25").Value = 9
$ws.Range("AH25").Value = 17
$ws.Range("AI25").Value = 17

# Row 39
$ws.Range("G39").Value = 2.2
$ws.Range("H39").Value = 3
$ws.Range("I39").Value = 3.5
$ws.Range("R39").Value = 1.8
$ws.Range("S39").Value = 1.91
$ws.Range("T39").Value = 7.5
$ws.Range("U39").Value = 10
$ws.Range("V39").Value = 9.5
$ws.Range("W39").Value = 21
$ws.Range("X39").Value = 19
$ws.Range("AB39").Value = 13
$ws.Range("AE39").Value = 10
$ws.Range("AF39").Value = 17
$ws.Range("AI39").Value = 29
